# Split several run-text blocks into multiple <w:t> runs joined by <w:br/>
# manual line breaks, mirroring the paragraph structure from the source diff.
$d = $word.ActiveDocument

# Objetivos (PT)
$findText = "Transmitir aos alunos os conceitos básicos relacionados diretamente a engenharia de sistemas biológicos capacitando-os ao entendimento dos princípios de engenharia envolvidos em operações em larga escala, em sistemas com organismos vivos, ecossistemas e processos biológicos.Apresentar aos alunos uma visão das aplicações potenciais e estratégicas da biotecnologia moderna, permitindo aos alunos estudar tópicos avançados em Engenharia de biossistemas, em uma abordagemvariável e multidisciplinar em temas relevantes a Engenharia.Aprimorar o raciocínio e despertar o espírito crítico e a criatividade dos alunos"
$replaceText = "Transmitir aos alunos os conceitos básicos relacionados diretamente a engenharia de sistemas biológicos capacitando-os ao entendimento dos princípios de engenharia envolvidos em operações em larga escala, em sistemas com organismos vivos, ecossistemas e processos biológicos.^lApresentar aos alunos uma visão das aplicações potenciais e estratégicas da biotecnologia moderna, permitindo aos alunos estudar tópicos avançados em Engenharia de biossistemas, em uma abordagem^lvariável e multidisciplinar em temas relevantes a Engenharia.^lAprimorar o raciocínio e despertar o espírito crítico e a criatividade dos alunos"
$d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null

# Objetivos (EN, italic)
$findText = "Transmit to students the basic concepts directly related to the engineering of biological systems, enabling them to understand the engineering principles involved in large-scale operations, in systems with living organisms, ecosystems and biological processes.Present students with an insight into the potential and strategic applications of modern biotechnology, allowing students to study advanced topics in Biosystems Engineering, in variable and multidisciplinar approach in topics relevant to Engineering.Improve reasoning and awaken students’ critical spirit and creativity."
$replaceText = "Transmit to students the basic concepts directly related to the engineering of biological systems, enabling them to understand the engineering principles involved in large-scale operations, in systems with living organisms, ecosystems and biological processes.^lPresent students with an insight into the potential and strategic applications of modern biotechnology, allowing students to study advanced topics in Biosystems Engineering, in variable and multidisciplinar approach in topics relevant to Engineering.^lImprove reasoning and awaken students’ critical spirit and creativity."
$d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null

# Programa resumido (PT)
$findText = "Introdução. Aspectos de Engenharia em processos fermentativos e enzimáticos envolvendo sistemas biológicos. Processos biotecnológicos de importância industrial. Fundamentos de engenharia aplicados aos biossistemas de importância na agricultura, medicina, biotecnologia, biofármacos, bioprocessamento industrial e conservação ambiental.2 Análise de critérios de ampliação de escala em processos envolvendo sistemas biológicos. Introdução às técnicas de separação/purificação de produtos biotecnológicos."
$replaceText = "Introdução. Aspectos de Engenharia em processos fermentativos e enzimáticos envolvendo sistemas biológicos. Processos biotecnológicos de importância industrial. Fundamentos de engenharia aplicados aos biossistemas de importância na agricultura, medicina, biotecnologia, biofármacos, bioprocessamento industrial e conservação ambiental.^l2 Análise de critérios de ampliação de escala em processos envolvendo sistemas biológicos. Introdução às técnicas de separação/purificação de produtos biotecnológicos."
$d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null

# Programa resumido (EN, italic)
$findText = "Introduction. Engineering aspects in fermentative and enzymatic processes involving biological systems. Biotechnological processes of industrial importance. Engineering fundamentals applied to biosystems of importance in agriculture, medicine, biotechnology, biopharmaceuticals, industrial bioprocessing and environmental conservation.Analysis of scale up criteria in processes involving biological systems. Introduction to separation/purification techniques for biotechnological products."
$replaceText = "Introduction. Engineering aspects in fermentative and enzymatic processes involving biological systems. Biotechnological processes of industrial importance. Engineering fundamentals applied to biosystems of importance in agriculture, medicine, biotechnology, biopharmaceuticals, industrial bioprocessing and environmental conservation.^lAnalysis of scale up criteria in processes involving biological systems. Introduction to separation/purification techniques for biotechnological products."
$d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null

# Programa (PT)
$findText = "- Introdução: importância dos bioprocessos e biossistemas e aplicações industriais.- Aspectos de Engenharia aplicados em Processos fermentativos e enzimáticos: características, biorreatores, operações, controle, sensores utilizados, aspectos cinéticos e modelagem de biossistemas.- Processos biotecnológicos de importância industrial: descrição e estudo de casos de alguns processos biotecnológicos.- Fundamentos de engenharia de bioprocessos aplicados aos biossistemas utilizando organismos vivos: transferência de oxigênio e respiração microbiana: transferência de massa (transferência por convecção em sistema gás-líquido; respiração microbiana; transferência de O2 da bolha de gás para a célula); transferência de O2 em biorreator (efeitos dos aspectos do dimensionamento e operacionais do biorreator - bolhas, aeração, agitação e propriedades do meio, agentes antiespumantes, temperatura, pressão do gás e pressão parcial de oxigênio). Transferência de potência e oxigênio em biorreator agitado e aerado.- Análise de critérios de variação de escala em processos envolvendo sistemas biológicos.- Fundamentos de engenharia aplicados aos biossistemas de importância na agricultura, medicina, biotecnologia, biofármacos, bioprocessamento industrial e conservação ambiental, exemplos práticos e estudo de casos.- Introdução às técnicas de separação/purificação de produtos biotecnológicos."
$replaceText = "- Introdução: importância dos bioprocessos e biossistemas e aplicações industriais.^l- Aspectos de Engenharia aplicados em Processos fermentativos e enzimáticos: características, biorreatores, operações, controle, sensores utilizados, aspectos cinéticos e modelagem de biossistemas.^l- Processos biotecnológicos de importância industrial: descrição e estudo de casos de alguns processos biotecnológicos.^l- Fundamentos de engenharia de bioprocessos aplicados aos biossistemas utilizando organismos vivos: transferência de oxigênio e respiração microbiana: transferência de massa (transferência por convecção em sistema gás-líquido; respiração microbiana; transferência de O2 da bolha de gás para a célula); transferência de O2 em biorreator (efeitos dos aspectos do dimensionamento e operacionais do biorreator - bolhas, aeração, agitação e propriedades do meio, agentes antiespumantes, temperatura, pressão do gás e pressão parcial de oxigênio). Transferência de potência e oxigênio em biorreator agitado e aerado.^l- Análise de critérios de variação de escala em processos envolvendo sistemas biológicos.^l- Fundamentos de engenharia aplicados aos biossistemas de importância na agricultura, medicina, biotecnologia, biofármacos, bioprocessamento industrial e conservação ambiental, exemplos práticos e estudo de casos.^l- Introdução às técnicas de separação/purificação de produtos biotecnológicos."
$d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null

# Programa (EN, italic)
$findText = "- Introduction: importance of bioprocesses and biosystems and industrial applications- Engineering aspects applied to fermentative and enzymatic processes: characteristics, bioreactors, operations, control, sensors used, kinetic aspects and modeling of biosystems.- Biotechnological processes of industrial importance: description and case study of some biotechnological processes- Fundamentals of bioprocess engineering applied to biosystems using living organisms: oxygen transfer and microbial respiration: mass transfer (convection transfer in a gas-liquid system; microbial respiration; transfer of O2 from the gas bubble to the cell); O2 transfer in bioreactor (effects of sizing and operational aspects of the bioreactor - bubbles, aeration, agitation and3 medium properties, antifoaming agents, temperature, gas pressure and partial pressure ofoxygen). Power and oxygen transfer in a stirred and aerated bioreactor.- Analysis of scale up variation criteria in processes involving biological systems.- Engineering fundamentals applied to biosystems of importance in agriculture, medicine, biotechnology, biopharmaceuticals, industrial bioprocessing and environmental conservation, practical examples and case studies.- Introduction to separation/purification techniques for biotechnological products."
$replaceText = "- Introduction: importance of bioprocesses and biosystems and industrial applications^l- Engineering aspects applied to fermentative and enzymatic processes: characteristics, bioreactors, operations, control, sensors used, kinetic aspects and modeling of biosystems.^l- Biotechnological processes of industrial importance: description and case study of some biotechnological processes^l- Fundamentals of bioprocess engineering applied to biosystems using living organisms: oxygen transfer and microbial respiration: mass transfer (convection transfer in a gas-liquid system; microbial respiration; transfer of O2 from the gas bubble to the cell); O2 transfer in bioreactor (effects of sizing and operational aspects of the bioreactor - bubbles, aeration, agitation and^l3 medium properties, antifoaming agents, temperature, gas pressure and partial pressure of^loxygen). Power and oxygen transfer in a stirred and aerated bioreactor.^l- Analysis of scale up variation criteria in processes involving biological systems.^l- Engineering fundamentals applied to biosystems of importance in agriculture, medicine, biotechnology, biopharmaceuticals, industrial bioprocessing and environmental conservation, practical examples and case studies.^l- Introduction to separation/purification techniques for biotechnological products."
$d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null

# Bibliografia
$findText = "1. LIMA, U.A. et al. Biotecnologia Industrial, vol. 3 - Processos Fermentativos e Enzimáticos - 1ª ed. - Edgard Blucher, 20012. LIMA, U.A. et al. Biotecnologia Industrial, vol. 3 - Processos Fermentativos e Enzimáticos - 2ª ed. - Edgard Blucher, 2020.3. AQUARONE, E. et al. Biotecnologia Industrial, vol. 4 – Biotecnologia na Produção de Alimentos - Edgard Blucher, 2001. 4. CASTILHO, L.R.; AUGUSTO, E.F.P.; MORAES, A. Tecnologia de Cultivo de Células Animais - de Biofármacos à Terapia Gênica. Roca, 2008.5. PESSOA JR, Adalberto et al. Biotecnologia farmacêutica: Aspectos sobre aplicação industrial. Editora Blucher, 2021.6. NASCIMENTO, R. et al. Microbiologia Industrial, vol. 1. - Bioprocessos. Elsevier, 2017.7. NASCIMENTO, R. et al. Microbiologia Industrial, vol. 2 - Alimentos. Elsevier, 2017."
$replaceText = "1. LIMA, U.A. et al. Biotecnologia Industrial, vol. 3 - Processos Fermentativos e Enzimáticos - 1ª ed. - Edgard Blucher, 2001^l2. LIMA, U.A. et al. Biotecnologia Industrial, vol. 3 - Processos Fermentativos e Enzimáticos - 2ª ed. - Edgard Blucher, 2020.^l3. AQUARONE, E. et al. Biotecnologia Industrial, vol. 4 – Biotecnologia na Produção de Alimentos - Edgard Blucher, 2001. ^l4. CASTILHO, L.R.; AUGUSTO, E.F.P.; MORAES, A. Tecnologia de Cultivo de Células Animais - de Biofármacos à Terapia Gênica. Roca, 2008.^l5. PESSOA JR, Adalberto et al. Biotecnologia farmacêutica: Aspectos sobre aplicação industrial. Editora Blucher, 2021.^l6. NASCIMENTO, R. et al. Microbiologia Industrial, vol. 1. - Bioprocessos. Elsevier, 2017.^l7. NASCIMENTO, R. et al. Microbiologia Industrial, vol. 2 - Alimentos. Elsevier, 2017."
$d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
